$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 6, shifting existing rows 6-7 down to 7-8.
$ws.Rows.Item(6).Insert()

# Populate the new row 6 with this week's data (weekly Fruta/Hortaliza update).
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C6").Value = "Arica y Parinacota"
$ws.Range("D6").Value = 44489
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = "Fruta"
$ws.Range("G6").Value = 100101
$ws.Range("H6").Value = "Berries"
$ws.Range("I6").Value = 100101007
$ws.Range("J6").Value = "Kiwi"
$ws.Range("K6").Value = "Hayward"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 26000
$ws.Range("O6").Value = 27000
$ws.Range("P6").Value = 26500
$ws.Range("Q6").Value = "`$/bandeja 18 kilos"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1472
$ws.Range("T6").Value = 18
